$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "Login using credentials from Excel with Apache POI"
$ws.Range("B19").Value = "FAILED"
$ws.Range("C19").Value = "edge"

$ws.Range("A20").Value = "Login using credentials from Excel with Apache POI"
$ws.Range("B20").Value = "PASSED"
$ws.Range("C20").Value = "edge"
